$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move existing row 3 (sairj.pdf) down to row 5 first, so it isn't clobbered
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "sairj.pdf"
$ws.Range("C5").Value = 22
$ws.Range("D5").Value = 24.4
$ws.Range("E5").Value = 23.2
$ws.Range("F5").Value = "python, statistics, pandas, data analysis"

# Move existing row 2 (sar.pdf) down to row 3
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "sar.pdf"
$ws.Range("C3").Value = 64
$ws.Range("D3").Value = 60.8
$ws.Range("E3").Value = 62.4
$ws.Range("F3").Value = "statistics"

# New row 2: Shardul mode.pdf
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Shardul mode.pdf"
$ws.Range("C2").Value = 70
$ws.Range("D2").Value = 66
$ws.Range("E2").Value = 68
$ws.Range("F2").Value = "statistics"

# New row 4: sairaj pawar.pdf
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "sairaj pawar.pdf"
$ws.Range("C4").Value = 64
$ws.Range("D4").Value = 60.8
$ws.Range("E4").Value = 62.4
$ws.Range("F4").Value = "statistics"

# New row 6: subhash dev.pdf
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "subhash dev.pdf"
$ws.Range("C6").Value = 10
$ws.Range("D6").Value = 14
$ws.Range("E6").Value = 12
$ws.Range("F6").Value = "python, machine learning, statistics, pandas, data analysis"
